# Update countries & provincias Spain
# - Refresh case numbers for several provinces (Madrid, Cataluna, Castilla y Leon,
#   Castilla-La Mancha, Pais Vasco, Aragon, Araba/Alava, La Rioja, Asturias, Gran Canaria)
# - Aragon's "Casos totales" overtakes Valencia/Valencia's, so the two rows swap order
# - Bump the "Datos actualizados" timestamp from 17:35 to 18:05

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 / 14: Aragon now outranks Valencia/Valencia, so they swap places ---
$ws.Range("A13").Value = "Aragon"
$ws.Range("B13").Value = 5618
$ws.Range("C13").Value = 3772
$ws.Range("D13").Value = 4760
$ws.Range("E13").Value = 858

$ws.Range("A14").Value = "Valencia/Valencia"
$ws.Range("B14").Value = 5609
$ws.Range("C14").Value = 4907
$ws.Range("D14").Value = 2767
$ws.Range("E14").Value = 693

# --- Madrid (row 4) ---
$ws.Range("B4").Value = 67425
$ws.Range("D4").Value = 58481
$ws.Range("E4").Value = 8944

# --- Cataluna (row 5) ---
$ws.Range("B5").Value = 57036
$ws.Range("D5").Value = 50380
$ws.Range("E5").Value = 6656

# --- Castilla y Leon (row 6) ---
$ws.Range("B6").Value = 18674
$ws.Range("D6").Value = 16713
$ws.Range("E6").Value = 1961

# --- Castilla-La Mancha (row 7) ---
$ws.Range("B7").Value = 16830
$ws.Range("D7").Value = 13901
$ws.Range("E7").Value = 2929

# --- Pais Vasco (row 9) ---
$ws.Range("B9").Value = 12576
$ws.Range("D9").Value = 11199
$ws.Range("E9").Value = 1377

# --- Araba/Alava (row 16) ---
$ws.Range("B16").Value = 5200
$ws.Range("D16").Value = 4692
$ws.Range("E16").Value = 508

# --- La Rioja (row 20) ---
$ws.Range("B20").Value = 4036
$ws.Range("D20").Value = 3680
$ws.Range("E20").Value = 356

# --- Asturias (row 32) ---
$ws.Range("B32").Value = 2376
$ws.Range("D32").Value = 2069

# --- Gran Canaria (row 33) ---
$ws.Range("B33").Value = 2312
$ws.Range("D33").Value = 2157

# --- Timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 18:05"
